$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the existing B5 timestamp value (tiny floating-point refinement)
$ws.Range("B5").Value = 43594.64684027381

# Add new row 6 data
$ws.Range("A6").Value = "Brevemente"

$ws.Range("B6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B6").Value = 44145.65398057337

$ws.Range("C6").Value = "dadadada"
$ws.Range("D6").Value = "Produtos para Gatos"
$ws.Range("E6").Value = 3

# F6 must stay a text string "02" (not be auto-converted to the number 2)
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "02"

$ws.Range("G6").Value = "Fabricante: ROYAL CANIN AROMATIC EXIGENT; "
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 4.78
$ws.Range("J6").Value = 4.78
